# Updated cryptos list on Tue Jun 20 23:18:51 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for
# every coin row on the sheet, and fixes the ranking order of HuobiToken
# and ImmutableX (rows 35-36), whose relative rank flipped in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.192.97'
$ws.Range("E2").Value = '  +5.52%  '
$ws.Range("D3").Value = '1.787.34'
$ws.Range("E3").Value = '  +3.16%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2689'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06291'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("D10").Value = '1.783.03'
$ws.Range("E10").Value = '  +2.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07052'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6288'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.664'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '80.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("D16").Value = '28.172.23'
$ws.Range("E16").Value = '  +6.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007240'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.59%  '
$ws.Range("D21").Value = '2.010.18'
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.556'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.760'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.249'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.853'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '109.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.180'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.19%  '
$ws.Range("E31").Value = '  +3.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.772'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04895'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.086'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.34%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6559'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.21%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.619'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9447'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.617'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.057'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.922'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01548'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.201'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1216'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05449'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.054'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.297'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.40%  '
